# Properties file will control which users to run smoke tests for.
#
# The "ignore" column (C) on Sheet1 flags whether a row of test
# credentials should be skipped. Flip the NU (Noor.Uddin) rows to be
# ignored and enable the PG (priya.giri) rows instead, adding the two
# missing Manufacturer/AuthorisedRep rows for PG that mirror the
# existing NU ones.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Flip NU rows (5-7) from "yes" to "no" (now ignored).
$ws.Range("C5").Value = "no"
$ws.Range("C6").Value = "no"
$ws.Range("C7").Value = "no"

# Flip PG rows (20-22) from "no" to "yes" (now active).
$ws.Range("C20").Value = "yes"
$ws.Range("C21").Value = "yes"
$ws.Range("C22").Value = "yes"

# Replace the duplicated NU username placeholders on rows 21/22 with the
# proper PG-specific usernames.
$ws.Range("A21").Value = "Manufacturer78H14_PG"
$ws.Range("A22").Value = "AuthorisedRep78H14_PG"

# Move the active selection to C7.
$ws.Range("C7").Select()
